$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the per-field rows (name / cost / type / text...) into a single
# Python-tuple-style string per card, one row per card.

$ws.Range("A2").Value = "('Corrupt', ['{5}{B}', 'Sorcery', 'Corrupt deals damage to any target equal to the number of Swamps you control. You gain life equal to the damage dealt this way.'])"
$ws.Range("A3").Value = "('Damnation', ['{2}{B}{B}', 'Sorcery', 'Destroy all creatures. They can" + [char]8217 + "t be regenerated.'])"
$ws.Range("A4").Value = "('Harmonize', ['{2}{G}{G}', 'Sorcery', 'Draw three cards.'])"
$ws.Range("A5").Value = "('Incinerate', ['{1}{R}', 'Instant', 'Incinerate deals 3 damage to any target. A creature dealt damage this way can" + [char]8217 + "t be regenerated this turn.'])"
$ws.Range("A6").Value = "('Mana Tithe', ['{W}', 'Instant', 'Counter target spell unless its controller pays {1}.'])"
$ws.Range("A7").Value = "('Ponder', ['{U}', 'Sorcery', 'Look at the top three cards of your library, then put them back in any order. You may shuffle your library.', 'Draw a card.'])"
$ws.Range("A8").Value = "('Tidings', ['{3}{U}{U}', 'Sorcery', 'Draw four cards.'])"

# Remove the now-unused rows 9 through 30
$ws.Range("A9:A30").ClearContents()
